$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999019174424,
    0.99888287058427372,
    0.99645592675115324,
    1.0005293919720049,
    0.99149140915699796,
    0.97172043236307193,
    0.97004772537702344,
    0.96715380890593938,
    0.9667046726576356,
    0.96720211046493032,
    0.96730321755216919,
    0.96775258872670966,
    0.97295545779516113,
    0.97650824436627848,
    0.974450051409818,
    0.97241152104201145,
    0.97086794467908932,
    0.97029049934398193,
    0.99782971717511293,
    0.99071295046317742,
    0.98931450446673419,
    0.98805001049493646,
    0.98942203062549572,
    0.97640205927429569,
    0.96994522400304772,
    0.9687840833218675,
    0.96493832827211456,
    0.95056820335931413,
    0.94054146476672762,
    0.93566735955278868,
    0.93529823296298154,
    0.93642516972771017,
    0.9401462086730803
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
